# Weekly update: two new daily price records are inserted into the
# "Fruta, Terminal Hortofrutícola Agro Chillán - Mango" data block.
# One new row is inserted right before the existing row 129, and a second
# new row is inserted right before the (shifted) last data row, pushing the
# previous final row down to become the very last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 129 (shifts old 129..159 down to 130..160) ---
$ws.Rows.Item(129).Insert()

$ws.Cells.Item(129, 1).Value = 7
$ws.Cells.Item(129, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(129, 3).Value = "Ñuble"
$ws.Cells.Item(129, 4).Value = 45120
$ws.Cells.Item(129, 5).Value = 16
$ws.Cells.Item(129, 6).Value = "Fruta"
$ws.Cells.Item(129, 7).Value = 100108
$ws.Cells.Item(129, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(129, 9).Value = 100108002
$ws.Cells.Item(129, 10).Value = "Mango"
$ws.Cells.Item(129, 11).Value = "Sin especificar"
$ws.Cells.Item(129, 12).Value = "Primera"
$ws.Cells.Item(129, 13).Value = 60
$ws.Cells.Item(129, 14).Value = 9000
$ws.Cells.Item(129, 15).Value = 9000
$ws.Cells.Item(129, 16).Value = 9000
$ws.Cells.Item(129, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(129, 18).Value = "Brasil"
$ws.Cells.Item(129, 19).Value = 2250
$ws.Cells.Item(129, 20).Value = 4

# --- Insert second new row at position 160 (shifts old final row 160 -> 161) ---
$ws.Rows.Item(160).Insert()

$ws.Cells.Item(160, 1).Value = 7
$ws.Cells.Item(160, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(160, 3).Value = "Ñuble"
$ws.Cells.Item(160, 4).Value = 45121
$ws.Cells.Item(160, 5).Value = 16
$ws.Cells.Item(160, 6).Value = "Fruta"
$ws.Cells.Item(160, 7).Value = 100108
$ws.Cells.Item(160, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(160, 9).Value = 100108002
$ws.Cells.Item(160, 10).Value = "Mango"
$ws.Cells.Item(160, 11).Value = "Sin especificar"
$ws.Cells.Item(160, 12).Value = "Primera"
$ws.Cells.Item(160, 13).Value = 40
$ws.Cells.Item(160, 14).Value = 9000
$ws.Cells.Item(160, 15).Value = 9000
$ws.Cells.Item(160, 16).Value = 9000
$ws.Cells.Item(160, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(160, 18).Value = "Brasil"
$ws.Cells.Item(160, 19).Value = 2250
$ws.Cells.Item(160, 20).Value = 4
